$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Countries & provincias Spain update (24 Apr 2020, 08:22 -> 08:52 refresh)
#
# The source feed re-sorted several neighbouring rows (by total cases) and
# refreshed their case counts. Net effect on this sheet, row-by-row:
#   Row 41: now Ucrania (new case numbers)
#   Row 42: now Noruega (row 41's old numbers, country shifted down one row)
#   Row 43: now Serbia  (row 42's old numbers, country shifted down one row)
#   Row 44: now Chequia (refreshed case numbers)
#   Row 109: now Georgia (new case numbers)
#   Row 110: now Taiwan  (row 109's old numbers, country shifted down one row)
#   Row 171: same country (San Martin (Parte Francesa)), refreshed D/E only
#   Row 193: now San Vicente y las Granadinas (new case numbers)
#   Row 194: now Curazao (row 193's old numbers, country shifted down one row)
#   Row 210: now Bonaire, San Eustaquio y Saba (new case numbers)
#   Row 211: now Sudan del Sur (row 210's old numbers, country shifted down one row)
# ---------------------------------------------------------------------------

function Set-CountryRow {
    param($Row, $Country, $Total, $New, $Active, $Recovered, $Critical, $DeathsToday, $Deaths)
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $New
    $ws.Cells.Item($Row, 4).Value = $Active
    $ws.Cells.Item($Row, 5).Value = $Recovered
    $ws.Cells.Item($Row, 6).Value = $Critical
    $ws.Cells.Item($Row, 7).Value = $DeathsToday
    $ws.Cells.Item($Row, 8).Value = $Deaths
}

Set-CountryRow 41 "Ucrania" 7647 477 601 6853 45 6 193
Set-CountryRow 42 "Noruega" 7401 0   32  7175 50 0 194
Set-CountryRow 43 "Serbia"  7276 0   1063 6074 96 0 139
Set-CountryRow 44 "Chequia" 7188 1   2186 4789 76 3 213

Set-CountryRow 109 "Georgia" 431 6 114 312 6 0 5
Set-CountryRow 110 "Taiwan"  428 1 264 158 0 0 6

$ws.Cells.Item(171, 4).Value = 20
$ws.Cells.Item(171, 5).Value = 16

Set-CountryRow 193 "San Vicente y las Granadinas" 14 1 5  9 0 0 0
Set-CountryRow 194 "Curazao"                       14 0 11 2 0 0 1

Set-CountryRow 210 "Bonaire, San Eustaquio y Saba" 5 0 0 5 0 0 0
Set-CountryRow 211 "Sudan del Sur"                 5 1 0 5 0 0 0

# Update the "last refreshed" footer string (row 1, last row of the sheet).
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 08:52"
